$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.681.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.65%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.726.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4919"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2622"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06227"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.02%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.85"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06984"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6120"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.505"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9981"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.32%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.509.53"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9986"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007202"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.43"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.950.69"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.462"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.78%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.571"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.111"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.11"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.79%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.389"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.914"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07984"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.673"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.46%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.77%  "

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9980"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.26%  "

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.609"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.30%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6271"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9323"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.95%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.032"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.43%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.416"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.36%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9997"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.28%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01514"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.58%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.566"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.50"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.97%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3862"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.922"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.86%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1157"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.20%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05382"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.25%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.884"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.88%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.37"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.17%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.76"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.16%  "
